$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("questions")

# Fix the wording of the "how do you know Amazon Connect" question.
$ws.Range("A3").Value = "How do you know Amazon Connect? A. News Letter, B. Social Media, C. AWS Event, D. AWS Website, or E. From Friend."

# Rename the callback-time question and give it an explicit black font color.
$ws.Range("A6").Value = "Preferred call back time?"
$ws.Range("A6").Font.Color = 0

# Leave the cursor on the last-edited cell of the questions sheet.
[void]$ws.Range("A6").Select()

# Switch to the receivers sheet and select B3, matching the saved selection state.
$ws3 = $wb.Worksheets.Item("receivers")
[void]$ws3.Activate()
[void]$ws3.Range("B3").Select()
